# Updates league odds bases: several match rows had their row order corrected.
# For each pair of adjacent rows below, every column except "A" (the sequential
# row counter) is swapped between the two rows, effectively swapping which
# match record occupies which row while keeping the A-column index sequence
# intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(9, 10),
    @(36, 37),
    @(49, 50),
    @(76, 77),
    @(122, 123),
    @(177, 178)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $addr1 = "B" + $r1 + ":AB" + $r1
    $addr2 = "B" + $r2 + ":AB" + $r2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
